$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.348.98"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "3.419.49"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.25"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.46"
$ws.Range("E6").Value = "  +2.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +4.71%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "3.417.75"
$ws.Range("E9").Value = "  -0.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.96"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").Value = "4.010.77"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.17"
$ws.Range("E15").Value = "  -2.54%  "

$ws.Range("D16").Value = "66.320.65"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").Value = "3.443.93"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.90"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.80"
$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.58"
$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.57"
$ws.Range("E22").Value = "  -2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.16"
$ws.Range("E23").Value = "  +1.82%  "

$ws.Range("E24").Value = "  +7.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.534"
$ws.Range("E26").Value = "  +0.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.80"
$ws.Range("E27").Value = "  +1.88%  "

$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.77"
$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.99"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.36"
$ws.Range("E32").Value = "  -2.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.01"
$ws.Range("E34").Value = "  -0.94%  "

$ws.Range("E35").Value = "  -4.12%  "

$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.70"
$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("E38").Value = "  -2.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.57"
$ws.Range("E39").Value = "  -5.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("E40").Value = "  +2.04%  "

$ws.Range("E41").Value = "  -3.33%  "

$ws.Range("E42").Value = "  +0.27%  "

$ws.Range("D43").Value = "2.714.61"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.30"
$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0686"
$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.89"
$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.02"
$ws.Range("E47").Value = "  -1.10%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "337.00"
$ws.Range("E48").Value = "  +9.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0287"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.89"
$ws.Range("E51").Value = "  +5.41%  "
